$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 32: fill in previously-empty task row ---
$ws.Range("A32").Value = "Ret SD0804 angivStraksAfskrivning"
$ws.Range("B32").Value = "Software Architect"
$ws.Range("C32").Value = 43894
$ws.Range("D32").Value = 0.375
$ws.Range("E32").Value = 0.47916666666666669

# --- Row 33: new task row ---
$ws.Range("A33").Value = "Lav DCD0804 angivStraksAfskrivning"
$ws.Range("B33").Value = "Software Architect"
$ws.Range("C33").Value = 43894
$ws.Range("D33").Value = 0.5
$ws.Range("E33").Value = 0.52083333333333337

# --- Row 34: new task row ---
$ws.Range("A34").Value = "Ret SD0804 angivStraksAfskrivning"
$ws.Range("B34").Value = "Software Architect"
$ws.Range("C34").Value = 43894
$ws.Range("D34").Value = 0.54166666666666663
$ws.Range("E34").Value = 0.59375

# --- Row 35: new task row ---
$ws.Range("A35").Value = "Kundemøde med HØK om UC10, DOM10, ATD10"
$ws.Range("B35").Value = "Requirement Specifier"
$ws.Range("C35").Value = 43894
$ws.Range("D35").Value = 0.59375
$ws.Range("E35").Value = 0.66666666666666663

# --- Extend the "time spent" shared formula down through row 35 ---
$ws.Range("G4:G35").Formula = "=E4-D4"

# --- Update the remembered selection to match the saved workbook view ---
$ws.Range("E36").Select()
